$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

$changes = [ordered]@{
    "D2" = "300.67"
    "E2" = "-0.71%"
    "D3" = "37.57"
    "E3" = "7.85%"
    "D4" = "5.001"
    "E4" = "-3.16%"
    "D5" = "0.07852"
    "E5" = "1.32%"
    "D6" = "2.199"
    "E6" = "-7.57%"
    "D7" = "8.006"
    "E7" = "-0.06%"
    "D8" = "3.995"
    "E8" = "1.18%"
    "D9" = "0.9120"
    "E9" = "-1.89%"
    "B10" = "WazirX"
    "C10" = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
    "D10" = "0.1865"
    "E10" = "3.82%"
    "B11" = "LiechtensteinCryptoassetsExchange"
    "C11" = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
    "D11" = "0.09227"
    "E11" = "-5.71%"
    "D12" = "0.08472"
    "E12" = "-2.18%"
    "D13" = "0.03521"
    "E13" = "6.02%"
    "D14" = "0.09923"
    "E14" = "0.26%"
    "D15" = "0.001476"
    "E15" = "-1.61%"
    "D16" = "0.005629"
    "E16" = "-2.36%"
    "D17" = "3.475"
    "E17" = "0.45%"
    "D18" = "2.095"
    "E18" = "-1.85%"
    "E19" = "2.84%"
    "D21" = "4.535"
    "E21" = "4.52%"
    "D22" = "0.2224"
    "E22" = "-3.34%"
    "D23" = "0.04637"
    "E23" = "1.31%"
    "E24" = "0.81%"
    "D25" = "0.004446"
    "E25" = "-0.35%"
    "E26" = "-0.34%"
    "D39" = "0.01749"
    "E39" = "-2.01%"
    "D40" = "0.04699"
    "E40" = "-2.06%"
    "D41" = "0.007851"
    "E41" = "1.27%"
    "D42" = "0.1388"
    "E42" = "-1.75%"
    "D43" = "0.007648"
    "E43" = "7.66%"
    "D44" = "0.002293"
    "E44" = "7.10%"
    "D45" = "0.01020"
    "E45" = "11.06%"
    "D46" = "0.00006048"
    "E46" = "-1.22%"
    "E47" = "-0.37%"
    "D48" = "8.676"
    "E48" = "183.26%"
    "E49" = "34.85%"
    "D50" = "0.00002094"
    "E50" = "-0.37%"
    "E51" = "-0.37%"
}

foreach ($cellRef in $changes.Keys) {
    Set-TextValue $cellRef $changes[$cellRef]
}
